$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44323
$ws.Range("M2").Value = 60

# Row 3
$ws.Range("D3").Value = 44323
$ws.Range("L3").Value = 'Segunda'
$ws.Range("M3").Value = 50

# Row 4
$ws.Range("D4").Value = 44322
$ws.Range("L4").Value = 'Primera'
$ws.Range("N4").Value = 10000
$ws.Range("O4").Value = 10000
$ws.Range("P4").Value = 10000
$ws.Range("S4").Value = 1000

# Row 5
$ws.Range("D5").Value = 44322
$ws.Range("L5").Value = 'Segunda'
$ws.Range("M5").Value = 40
$ws.Range("N5").Value = 8000
$ws.Range("O5").Value = 8000
$ws.Range("P5").Value = 8000
$ws.Range("S5").Value = 800

# Row 6
$ws.Range("D6").Value = 44328
$ws.Range("N6").Value = 8000
$ws.Range("O6").Value = 8000
$ws.Range("P6").Value = 8000
$ws.Range("S6").Value = 800

# Row 7
$ws.Range("D7").Value = 44328
$ws.Range("L7").Value = 'Segunda'
$ws.Range("M7").Value = 48
$ws.Range("N7").Value = 7000
$ws.Range("O7").Value = 7000
$ws.Range("P7").Value = 7000
$ws.Range("S7").Value = 700

# Row 8
$ws.Range("D8").Value = 44307
$ws.Range("L8").Value = 'Primera'
$ws.Range("M8").Value = 40
$ws.Range("N8").Value = 10000
$ws.Range("O8").Value = 10000
$ws.Range("P8").Value = 10000
$ws.Range("S8").Value = 1000

# Row 9
$ws.Range("D9").Value = 45082
$ws.Range("L9").Value = 'Especial'
$ws.Range("M9").Value = 56
$ws.Range("N9").Value = 15000
$ws.Range("O9").Value = 15000
$ws.Range("P9").Value = 15000
$ws.Range("R9").Value = 'Región de O''Higgins'
$ws.Range("S9").Value = 1500

# Row 10
$ws.Range("D10").Value = 45082
$ws.Range("L10").Value = 'Primera'
$ws.Range("N10").Value = 12000
$ws.Range("O10").Value = 12000
$ws.Range("P10").Value = 12000
$ws.Range("R10").Value = 'Región de O''Higgins'
$ws.Range("S10").Value = 1200

# Row 11
$ws.Range("D11").Value = 45082
$ws.Range("L11").Value = 'Segunda'
$ws.Range("M11").Value = 60
$ws.Range("R11").Value = 'Región de O''Higgins'

# Row 12
$ws.Range("D12").Value = 44699
$ws.Range("L12").Value = 'Especial'
$ws.Range("M12").Value = 56
$ws.Range("N12").Value = 12000
$ws.Range("O12").Value = 12000
$ws.Range("P12").Value = 12000
$ws.Range("S12").Value = 1200

# Row 13
$ws.Range("D13").Value = 44699
$ws.Range("M13").Value = 60

# Row 14
$ws.Range("D14").Value = 44309

# Row 15
$ws.Range("D15").Value = 44302
$ws.Range("M15").Value = 45
$ws.Range("N15").Value = 10000
$ws.Range("O15").Value = 10000
$ws.Range("P15").Value = 10000
$ws.Range("R15").Value = 'Provincia de Quillota'
$ws.Range("S15").Value = 1000

# Row 16
$ws.Range("D16").Value = 44312
$ws.Range("L16").Value = 'Primera'
$ws.Range("M16").Value = 48
$ws.Range("N16").Value = 10000
$ws.Range("O16").Value = 10000
$ws.Range("P16").Value = 10000
$ws.Range("R16").Value = 'Provincia de Quillota'
$ws.Range("S16").Value = 1000

# Row 17
$ws.Range("D17").Value = 44315
$ws.Range("M17").Value = 45

# Row 18
$ws.Range("D18").Value = 44329
$ws.Range("N18").Value = 9000
$ws.Range("O18").Value = 9000
$ws.Range("P18").Value = 9000
$ws.Range("R18").Value = 'Región Metropolitana'
$ws.Range("S18").Value = 900

# Row 19
$ws.Range("D19").Value = 44329
$ws.Range("M19").Value = 50
$ws.Range("R19").Value = 'Región Metropolitana'

# Row 20
$ws.Range("D20").Value = 44314
$ws.Range("L20").Value = 'Primera'
$ws.Range("M20").Value = 47
$ws.Range("N20").Value = 9000
$ws.Range("O20").Value = 9000
$ws.Range("P20").Value = 9000
$ws.Range("R20").Value = 'Provincia de Quillota'
$ws.Range("S20").Value = 900

# Row 21
$ws.Range("D21").Value = 44326
$ws.Range("M21").Value = 65
$ws.Range("N21").Value = 10000
$ws.Range("O21").Value = 10000
$ws.Range("P21").Value = 10000
$ws.Range("R21").Value = 'Provincia de Quillota'
$ws.Range("S21").Value = 1000

# Row 22
$ws.Range("D22").Value = 44326
$ws.Range("M22").Value = 67
$ws.Range("N22").Value = 8000
$ws.Range("O22").Value = 8000
$ws.Range("P22").Value = 8000
$ws.Range("R22").Value = 'Provincia de Quillota'
$ws.Range("S22").Value = 800

# Row 26
$ws.Range("D26").Value = 44321
$ws.Range("M26").Value = 58
$ws.Range("N26").Value = 9000
$ws.Range("O26").Value = 9000
$ws.Range("P26").Value = 9000
$ws.Range("S26").Value = 900

# Row 27
$ws.Range("D27").Value = 44333
$ws.Range("L27").Value = 'Especial'
$ws.Range("M27").Value = 58
$ws.Range("N27").Value = 10000
$ws.Range("O27").Value = 10000
$ws.Range("P27").Value = 10000
$ws.Range("S27").Value = 1000

# Row 28
$ws.Range("D28").Value = 44333
$ws.Range("M28").Value = 65
$ws.Range("N28").Value = 9000
$ws.Range("O28").Value = 9000
$ws.Range("P28").Value = 9000
$ws.Range("S28").Value = 900

# Row 29
$ws.Range("D29").Value = 44333
$ws.Range("M29").Value = 60

# Row 30
$ws.Range("D30").Value = 44301
$ws.Range("N30").Value = 10000
$ws.Range("O30").Value = 10000
$ws.Range("P30").Value = 10000
$ws.Range("S30").Value = 1000

# Row 31
$ws.Range("D31").Value = 44306
$ws.Range("L31").Value = 'Primera'
$ws.Range("M31").Value = 45
$ws.Range("N31").Value = 10000
$ws.Range("O31").Value = 10000
$ws.Range("P31").Value = 10000
$ws.Range("S31").Value = 1000

# Row 32
$ws.Range("D32").Value = 44308
$ws.Range("L32").Value = 'Primera'
$ws.Range("M32").Value = 45

# Row 33
$ws.Range("D33").Value = 44308
$ws.Range("L33").Value = 'Segunda'
$ws.Range("M33").Value = 48
$ws.Range("N33").Value = 8000
$ws.Range("O33").Value = 8000
$ws.Range("P33").Value = 8000
$ws.Range("S33").Value = 800

# Row 34
$ws.Range("D34").Value = 44319
$ws.Range("L34").Value = 'Primera'
$ws.Range("M34").Value = 68
$ws.Range("N34").Value = 10000
$ws.Range("O34").Value = 10000
$ws.Range("P34").Value = 10000
$ws.Range("S34").Value = 1000

# Row 35
$ws.Range("D35").Value = 44319
$ws.Range("M35").Value = 57
$ws.Range("N35").Value = 8000
$ws.Range("O35").Value = 8000
$ws.Range("P35").Value = 8000
$ws.Range("S35").Value = 800
